$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Renumber the Name_x / DeviceName_x columns (D/E) so they run sequentially
# 1..10 instead of having a gap at "Name_4"/"DeviceName_1"/"DeviceName_5".
$ws.Range("D2").Value = "Name_1"
$ws.Range("D3").Value = "Name_2"
$ws.Range("D4").Value = "Name_3"
$ws.Range("D5").Value = "Name_4"
$ws.Range("D6").Value = "Name_5"
$ws.Range("D7").Value = "Name_6"
$ws.Range("D8").Value = "Name_7"
$ws.Range("D9").Value = "Name_8"
$ws.Range("D10").Value = "Name_9"
$ws.Range("D11").Value = "Name_10"

$ws.Range("E2").Value = "DeviceName_1"
$ws.Range("E3").Value = "DeviceName_2"
$ws.Range("E4").Value = "DeviceName_3"
$ws.Range("E5").Value = "DeviceName_4"
$ws.Range("E6").Value = "DeviceName_5"
$ws.Range("E7").Value = "DeviceName_6"
$ws.Range("E8").Value = "DeviceName_7"
$ws.Range("E9").Value = "DeviceName_8"
$ws.Range("E10").Value = "DeviceName_9"
$ws.Range("E11").Value = "DeviceName_10"

# Update the view state: reset zoom to 100% and move the selection.
$excel.ActiveWindow.Zoom = 100
[void]$ws.Range("G20").Select()
